# Update the date heading
$d = $word.ActiveDocument
$d.Content.Find.Execute("2024-07-03 Wednesday", $true, $false, $false, $false, $false, $true, 1, $false, "2024-07-04 Thursday", 2) | Out-Null

# Update each answer cell in the table, in row-major order (20 rows x 5 cols)
$newValues = @(
    "76-8=68",
    "76+9=85",
    "30+63=93",
    "93-65=28",
    "33+8=41",
    "9-4=5",
    "31+42=73",
    "52+12=64",
    "13+72=85",
    "47+52=99",
    "51-49=2",
    "97-86=11",
    "1+7=8",
    "80-18=62",
    "8+39=47",
    "59-31=28",
    "39+39=78",
    "68-13=55",
    "16+19=35",
    "79+14=93",
    "14+46=60",
    "75-4=71",
    "56-46=10",
    "94-6=88",
    "30+46=76",
    "20+8=28",
    "14+5=19",
    "26+63=89",
    "59-49=10",
    "28+70=98",
    "19+37=56",
    "36-10=26",
    "25+10=35",
    "19+64=83",
    "12-10=2",
    "89-79=10",
    "37+19=56",
    "92-15=77",
    "33+34=67",
    "89-61=28",
    "95-78=17",
    "76-16=60",
    "51+2=53",
    "64-18=46",
    "43+9=52",
    "44-44=0",
    "10+65=75",
    "38-15=23",
    "9+83=92",
    "31-20=11",
    "50+24=74",
    "64+28=92",
    "45+35=80",
    "72+13=85",
    "73-13=60",
    "87-85=2",
    "45+45=90",
    "89-68=21",
    "30-27=3",
    "85-5=80",
    "61-4=57",
    "2+21=23",
    "33+43=76",
    "83-36=47",
    "48-27=21",
    "62+35=97",
    "89-47=42",
    "58+12=70",
    "24-10=14",
    "76-5=71",
    "10+62=72",
    "30+67=97",
    "33+40=73",
    "27+33=60",
    "24+52=76",
    "72-14=58",
    "87-79=8",
    "8+27=35",
    "72-22=50",
    "41-34=7",
    "99-60=39",
    "30-3=27",
    "77-16=61",
    "60-31=29",
    "4+50=54",
    "74-63=11",
    "38+60=98",
    "23+41=64",
    "37+27=64",
    "91-62=29",
    "19-11=8",
    "77+11=88",
    "84-63=21",
    "41+14=55",
    "87-32=55",
    "35+27=62",
    "71-47=24",
    "86-0=86",
    "60+30=90",
    "38+37=75"
)

$tbl = $d.Tables.Item(1)
$rows = $tbl.Rows.Count
$cols = $tbl.Columns.Count
$idx = 0
for ($r = 1; $r -le $rows; $r++) {
    for ($c = 1; $c -le $cols; $c++) {
        $cell = $tbl.Cell($r, $c)
        $cell.Range.Text = $newValues[$idx]
        $idx++
    }
}

Write-Output "Updated $idx cells"
